# Translate the "The ants problem - subtitles" document's English
# subtitle text to Swahili, run by run, using Find/Replace (wdReplaceAll).
# The very last replacement ("velocity" -> "kasi") is intentionally kept
# last: that single English word is also a substring of two longer
# sentences translated earlier in the script, so it must run only after
# those sentences have already been turned into Swahili.
$d = $word.ActiveDocument

$d.Content.Find.Execute('The ants problem - subtitles:', $true, $false, $false, $false, $false, $true, 1, $false, 'Tatizo la mchwa - manukuu:', 2) | Out-Null
$d.Content.Find.Execute('The dialogue starts at 40 seconds in so I added 27 seconds to the times as they were - John Argentino', $true, $false, $false, $false, $false, $true, 1, $false, 'Mazungumzo huanza kwa sekunde 40 kwa hivyo niliongeza sekunde 27 kwa nyakati kama zilivyokuwa - John Argentino', 2) | Out-Null
$d.Content.Find.Execute('[Music]', $true, $false, $false, $false, $false, $true, 1, $false, '[Muziki]', 2) | Out-Null
$d.Content.Find.Execute('okay so the puzzles I''m going to', $true, $false, $false, $false, $false, $true, 1, $false, 'sawa kwa hivyo mafumbo nitaenda', 2) | Out-Null
$d.Content.Find.Execute('challenge you with are two basic', $true, $false, $false, $false, $false, $true, 1, $false, 'changamoto uliyonayo ni mbili za msingi', 2) | Out-Null
$d.Content.Find.Execute('versions of a more complicated puzzle', $true, $false, $false, $false, $false, $true, 1, $false, 'matoleo ya fumbo ngumu zaidi', 2) | Out-Null
$d.Content.Find.Execute('known as the ants puzzle, which I''m', $true, $false, $false, $false, $false, $true, 1, $false, 'inayojulikana kama fumbo la mchwa, ambalo mimi ni', 2) | Out-Null
$d.Content.Find.Execute('probably going to discuss in a different', $true, $false, $false, $false, $false, $true, 1, $false, 'pengine kwenda kujadili katika tofauti', 2) | Out-Null
$d.Content.Find.Execute('video. Let me just finish writing down', $true, $false, $false, $false, $false, $true, 1, $false, 'video. Ngoja nimalizie kuandika', 2) | Out-Null
$d.Content.Find.Execute('the title and, well, I can even draw a', $true, $false, $false, $false, $false, $true, 1, $false, 'kichwa na, vizuri, naweza hata kuchora a', 2) | Out-Null
$d.Content.Find.Execute('little ant right here. okay, let''s get', $true, $false, $false, $false, $false, $true, 1, $false, 'mchwa mdogo hapa. sawa, tupate', 2) | Out-Null
$d.Content.Find.Execute('started! As I said I''m going to discuss', $true, $false, $false, $false, $false, $true, 1, $false, 'imeanza! Kama nilivyosema nitajadili', 2) | Out-Null
$d.Content.Find.Execute('two puzzles in the first puzzle there', $true, $false, $false, $false, $false, $true, 1, $false, 'mafumbo mawili katika fumbo la kwanza hapo', 2) | Out-Null
$d.Content.Find.Execute('are two ants on a very high stool: a sort', $true, $false, $false, $false, $false, $true, 1, $false, 'ni mchwa wawili kwenye kinyesi cha juu sana: aina', 2) | Out-Null
$d.Content.Find.Execute('of Mountain, flat on the top with two', $true, $false, $false, $false, $false, $true, 1, $false, 'ya Mlima, gorofa juu na mbili', 2) | Out-Null
$d.Content.Find.Execute('steep cliffs to both the sides. The flat', $true, $false, $false, $false, $false, $true, 1, $false, 'miamba mikali kwa pande zote mbili. Gorofa', 2) | Out-Null
$d.Content.Find.Execute('peak is one meter wide the two ants move', $true, $false, $false, $false, $false, $true, 1, $false, 'kilele ni mita moja upana wa mchwa wawili hoja', 2) | Out-Null
$d.Content.Find.Execute('with a velocity, let''s call it V, which is', $true, $false, $false, $false, $false, $true, 1, $false, 'kwa kasi, tuiite V, ambayo ni', 2) | Out-Null
$d.Content.Find.Execute('the same for both of them and that is', $true, $false, $false, $false, $false, $true, 1, $false, 'sawa kwa wote wawili na hiyo ni', 2) | Out-Null
$d.Content.Find.Execute('equal to one centimeter per second. You', $true, $false, $false, $false, $false, $true, 1, $false, 'sawa na sentimita moja kwa sekunde. Wewe', 2) | Out-Null
$d.Content.Find.Execute('can decide the direction towards each', $true, $false, $false, $false, $false, $true, 1, $false, 'inaweza kuamua mwelekeo kuelekea kila mmoja', 2) | Out-Null
$d.Content.Find.Execute('ant moves if it is right or left and', $true, $false, $false, $false, $false, $true, 1, $false, 'mchwa husogea ikiwa ni kulia au kushoto na', 2) | Out-Null
$d.Content.Find.Execute('where exactly to place the two ants on the', $true, $false, $false, $false, $false, $true, 1, $false, 'wapi hasa kuweka mchwa wawili kwenye', 2) | Out-Null
$d.Content.Find.Execute('top of the mountain. Your purpose is to', $true, $false, $false, $false, $false, $true, 1, $false, 'juu ya mlima. Kusudi lako ni', 2) | Out-Null
$d.Content.Find.Execute('make the time the last ant takes before', $true, $false, $false, $false, $false, $true, 1, $false, 'fanya wakati mchwa wa mwisho huchukua hapo awali', 2) | Out-Null
$d.Content.Find.Execute('falling the longest possible. Ants cannot', $true, $false, $false, $false, $false, $true, 1, $false, 'kuanguka kwa muda mrefu iwezekanavyo. Mchwa hawawezi', 2) | Out-Null
$d.Content.Find.Execute('be still: they must move to the right or', $true, $false, $false, $false, $false, $true, 1, $false, 'tulia: lazima wahamie kulia au', 2) | Out-Null
$d.Content.Find.Execute('to the left but they must move and after', $true, $false, $false, $false, $false, $true, 1, $false, 'upande wa kushoto lakini lazima wasogee na baada', 2) | Out-Null
$d.Content.Find.Execute('meeting each other they turn around and', $true, $false, $false, $false, $false, $true, 1, $false, 'wakikutana wanageuka na', 2) | Out-Null
$d.Content.Find.Execute('keep moving with the same but opposite', $true, $false, $false, $false, $false, $true, 1, $false, 'endelea kusonga na sawa lakini kinyume', 2) | Out-Null
$d.Content.Find.Execute('[Music]', $true, $false, $false, $false, $false, $true, 1, $false, '[Muziki]', 2) | Out-Null
$d.Content.Find.Execute('so again what are the precise positions', $true, $false, $false, $false, $false, $true, 1, $false, 'kwa hivyo tena ni nafasi gani sahihi', 2) | Out-Null
$d.Content.Find.Execute('where I should place the two ants in', $true, $false, $false, $false, $false, $true, 1, $false, 'ambapo ninapaswa kuwaweka mchwa wawili ndani', 2) | Out-Null
$d.Content.Find.Execute('order to get the longest time before the', $true, $false, $false, $false, $false, $true, 1, $false, 'ili kupata muda mrefu zaidi kabla ya', 2) | Out-Null
$d.Content.Find.Execute('last ant falls? The second puzzle is', $true, $false, $false, $false, $false, $true, 1, $false, 'chungu mwisho huanguka? Fumbo la pili ni', 2) | Out-Null
$d.Content.Find.Execute('basically the same but now we have three', $true, $false, $false, $false, $false, $true, 1, $false, 'kimsingi ni sawa lakini sasa tuna tatu', 2) | Out-Null
$d.Content.Find.Execute('ants instead of two.', $true, $false, $false, $false, $false, $true, 1, $false, 'mchwa badala ya wawili.', 2) | Out-Null
$d.Content.Find.Execute('As before the ants velocity is one', $true, $false, $false, $false, $false, $true, 1, $false, 'Kama kabla ya mchwa kasi ni moja', 2) | Out-Null
$d.Content.Find.Execute('centimeter per second, every ant turns', $true, $false, $false, $false, $false, $true, 1, $false, 'sentimita kwa sekunde, kila mchwa hugeuka', 2) | Out-Null
$d.Content.Find.Execute('around after meeting another ant and', $true, $false, $false, $false, $false, $true, 1, $false, 'karibu baada ya kukutana na mchwa mwingine na', 2) | Out-Null
$d.Content.Find.Execute('the peak is one meter wide. So, what are', $true, $false, $false, $false, $false, $true, 1, $false, 'kilele kina upana wa mita moja. Hivyo, ni nini', 2) | Out-Null
$d.Content.Find.Execute('now the precise positions', $true, $false, $false, $false, $false, $true, 1, $false, 'sasa nafasi sahihi', 2) | Out-Null
$d.Content.Find.Execute('I should place the three ants in order', $true, $false, $false, $false, $false, $true, 1, $false, 'Ninapaswa kuweka mchwa watatu kwa mpangilio', 2) | Out-Null
$d.Content.Find.Execute('to get the longest time before the last', $true, $false, $false, $false, $false, $true, 1, $false, 'kupata muda mrefu zaidi kabla ya mwisho', 2) | Out-Null
$d.Content.Find.Execute('ant falls down? I hope you enjoyed this', $true, $false, $false, $false, $false, $true, 1, $false, 'chungu huanguka chini? Natumaini ulifurahia hili', 2) | Out-Null
$d.Content.Find.Execute('video do your best and good luck', $true, $false, $false, $false, $false, $true, 1, $false, 'video fanya bora na bahati nzuri', 2) | Out-Null
$d.Content.Find.Execute('velocity', $true, $false, $false, $false, $false, $true, 1, $false, 'kasi', 2) | Out-Null
